$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestResults")

# Row 2 (existing row is reused/overwritten with new content)
$ws.Range("A2").Value = " iAU_TC_ID_205"
$ws.Range("B2").Value = "@RegressionA Pre-Request Verify Elumina Login and Create Exam"
$ws.Range("C2").Value = "passed"

# Row 3 (new)
$ws.Range("A3").Value = " iAU_TC_ID_205"
$ws.Range("B3").Value = '@RegressionA Pre-Request "Validation of Delivery --> Add New Users"'
$ws.Range("C3").Value = "passed"

# Row 4 (new)
$ws.Range("A4").Value = "iAU_TC_ID_180"
$ws.Range("B4").Value = "@RegressionA Validation of Manage Delivery --> Edit user "
$ws.Range("C4").Value = "passed"

# Row 5 (new)
$ws.Range("A5").Value = "iAU_TC_ID_205"
$ws.Range("B5").Value = "@RegressionA Validation of Delivery --> Venue Summary "
$ws.Range("C5").Value = "passed"

# Row 6 (new, no B6 value)
$ws.Range("A6").Value = "iAU_TC_ID_206.,iAU_TC_ID_210.,iAU_TC_ID_209.,iAU_TC_ID_211 @RegressionA Validation of Delivery --> Live Dashboard "
$ws.Range("C6").Value = "timedOut"
